$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (write A then B so shared-string order matches: Ticker, Mercado, ...)
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "Mercado"

# Column A: full list of tickers (written top-to-bottom before column B values,
# matching the order new strings were interned in the workbook)
$ws.Range("A2").Value = "BBAS3"
$ws.Range("A3").Value = "PETR4"
$ws.Range("A4").Value = "AMZO34"
$ws.Range("A5").Value = "HASH11"
$ws.Range("A6").Value = "CPTS11"
$ws.Range("A7").Value = "TESOURO IPCA+ 2026"

# Column B: market/segment for each ticker
$ws.Range("B2").Value = "Ações"
$ws.Range("B3").Value = "Ações"
$ws.Range("B4").Value = "BDR"
$ws.Range("B5").Value = "ETF"
$ws.Range("B6").Value = "FII"
$ws.Range("B7").Value = "Tesouro Direto"

# Best-fit style column widths (closest achievable snap to the original
# bestFit-computed widths of 19.5703125 / 14.28515625 characters)
$ws.Columns.Item(1).ColumnWidth = 18.65
$ws.Columns.Item(2).ColumnWidth = 13.5
